$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set row heights for rows 1-25 to 29
for ($r = 1; $r -le 25; $r++) {
    $ws.Rows.Item($r).RowHeight = 29
}

# Column widths for H,I,J,K (8,9,10,11)
$ws.Columns.Item(8).ColumnWidth = 39.1640625
$ws.Columns.Item(9).ColumnWidth = 39
$ws.Columns.Item(10).ColumnWidth = 36.33203125
$ws.Columns.Item(11).ColumnWidth = 44

# Selection / view
$ws.Range("G11").Select()
